$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 9: August (through 08-13) -> August (through 08-14), with updated values
$ws.Range("A9").Value = "August (through 08-14)"
$ws.Range("B9").Value = 14
$ws.Range("C9").Value = 32
$ws.Range("D9").Value = 31
$ws.Range("E9").Value = 23
$ws.Range("F9").Value = 19
$ws.Range("G9").Value = 86
$ws.Range("H9").Value = 77

# Row 10: Total values updated
$ws.Range("B10").Value = 176
$ws.Range("C10").Value = 334
$ws.Range("D10").Value = 496
$ws.Range("E10").Value = 448
$ws.Range("F10").Value = 323
$ws.Range("G10").Value = 707
$ws.Range("H10").Value = 992
